# "Add files via upload" - the workbook was re-saved/re-uploaded. The only
# substantive (non-cosmetic) change baked into the canonical OOXML is inside
# row 2 of Sheet1 / the shared-string table: cell B2 and cell C2 end up
# pointing at different <si> slots (B2: 3 -> 4, C2: 4 -> 3), while what they
# *display* is unchanged - B2 still holds the obfuscated/"encrypted" looking
# text (re-salted, since that kind of string is re-generated with a fresh
# random salt on every save) and C2 still holds "13:15-13:20".
#
# Reproduce that end state with the Excel object model: write the literal
# values back out in an order that makes the shared-string table land with
# "13:15-13:20" in the lower slot and the (newly salted) obfuscated text in
# the next slot - i.e. the same slot swap the diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$time = "13:15-13:20"
$obfuscated = "U2FsdGVkX1+IuIj1EkUoq5H/fWdtmWJJDlHwdRSVK/riF3fe3k79mZBMQsa2D3Cw4TBiZR5N327jbmal/0QTiSFryMdQh937Ohds9V0Tv7QUvA1LyNtmBhl9iEEFMV9wqGKMbg51r+SFA/J9C8i60wE51ZvuHoZQtaEfqNjCnkp8wQHkKDAZZE8AybBsZrsXGdjmx/G6oYBhi+KJSIzT+wy8AzPlJw47/sdcWugUfI0UIxpLqBhg5HJ9XF0PjV9iDNH4LdRRJZtRI9G5QT06Gr4fO9PgQrojRi83U8414QLHl55i7CGHR5JtLMgBh6z9hCqywrGYFpvxNU940kyefpqal3xtZ/tKDRlOvjc3znl+bVeeqJvx9HJGirw2e430MN1SPKYEr3nhdJbnkXriIpU31zcf3rU2PY5OpLXXtUEeCHQGI6+/qZhvhCudmbp54mDCTGfutjwYhls3aUImnTfhFa30/vpZ3L670m0IjxhBOajp8/DNkIwcwuVhp/fYyhHh0O9EzVje4ZmVQP8r/Gp8pF3nKtW9Sq6IiCc8p2kOYbv2yW3dX3dU+NMDIpfs9TMYamkqbwnaUf9m8qXBW/FCA3Hwaw8pPT83bQKnsIJ1LBzMgCUi9kyNg98GV5JqPXPG5WaWuhDh0fnMOc5Ry7Y3ys+E2V70HFRAJClKw9XktEUOKIK8bx8ITqMxTVg+v3K7C+lPrHqzQlP0sExaSZK596vvoZ1WGVWCSg6Z6JjWOV1RFFQYHfOWosNLc7j3xa/IP61CjxDRZOPF8K98iciS9LMldOTeMmqJzRmhjg/k8do3yDqt3yYnMID7B5+sbRBcRsDnR2BMptey580hczPa00tWwdYOpd/8DDX7j9wD3qmlPxbvUZRMDH03XpXxAnX8IfgZfeKQBIwnCBK8oVte5dObhF2gXX33S8o+l46ku9qfLokzQTyN0jt6hJw0qI0bafdopY3kEo9/MevGuLxHTWMGFXDBtNWscxrfNO4CbFILt2iO3kBThBA4cJO8x9T9DA9gkbUfKhFrJ2TM61DkYSDr1vVgIcrTSDphBiYhj4GYXMGZ9N5u5yX4L7FaFIEuyKy7odYK2vy5FyJnDlQvkq+GYdhgfPhWL6Fc7BeF0DN4Kn+flyDlERrYhu+EYDmzbpCTD/7dgp7ld2Vj/vKINAlpOlKuZ0p0gFoNAjFCN9ZmyQ1Yc6DSrdLf40fMfFrHv2MtgeFmeuvzN82Q4EZUA+/4+AjwsdLJEY2bX1bV3sg+o0lWNVFdECCdcW6Ib7tq/lzkrJrUgZXN2bOGRSjfvxdcBj8e1zRhj/nTtjM="

# Scratch cell: park the (unchanged) time string so it claims the lower
# shared-string slot before the obfuscated text is (re)written into B2.
$ws.Range("E1").Value = $time
$ws.Range("B2").Value = $obfuscated
$ws.Range("C2").Value = $time
$ws.Range("E1").ClearContents()
